$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows to reflect repulled data
$ws.Range("F2").Value = 1
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -3
$ws.Range("F18").Value = -3
$ws.Range("F19").Value = 3
$ws.Range("F21").Value = -6
$ws.Range("F22").Value = -1
$ws.Range("F24").Value = -5
$ws.Range("F25").Value = -6
$ws.Range("F26").Value = -1
